$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B rows 2-10: zero-padded ID codes, entered as text strings.
# Number format must be applied BEFORE the value is written so the "01"
# string is kept verbatim instead of being coerced to the number 1.
$ws.Range("B2:B10").NumberFormat = "@"
$ws.Range("B2").Value = "01"
$ws.Range("B3").Value = "02"
$ws.Range("B4").Value = "03"
$ws.Range("B5").Value = "04"
$ws.Range("B6").Value = "05"
$ws.Range("B7").Value = "06"
$ws.Range("B8").Value = "07"
$ws.Range("B9").Value = "08"
$ws.Range("B10").Value = "09"

# Column B rows 11-16 keep numeric values (10-15), with the text number
# format applied cosmetically AFTER the value is written (matches the
# target file, where these cells carry style s="1" but still store a
# numeric <v>).
$ws.Range("B11").Value = 10
$ws.Range("B12").Value = 11
$ws.Range("B13").Value = 12
$ws.Range("B14").Value = 13
$ws.Range("B15").Value = 14
$ws.Range("B16").Value = 15
$ws.Range("B11:B16").NumberFormat = "@"

# New column C header
$ws.Range("C1").Value = "Canopy"

# Column C: canopy classification for each land-cover type
$ws.Range("C2").Value = "Open"
$ws.Range("C3").Value = "Open"
$ws.Range("C4").Value = "Mixed"
$ws.Range("C5").Value = "Closed"
$ws.Range("C6").Value = "Closed"
$ws.Range("C7").Value = "Closed"
$ws.Range("C8").Value = "Open"
$ws.Range("C9").Value = "Open"
$ws.Range("C10").Value = "Mixed"
$ws.Range("C11").Value = "Open"
$ws.Range("C12").Value = "Open"
$ws.Range("C13").Value = "Open"
$ws.Range("C14").Value = "Mixed"
$ws.Range("C15").Value = "Mixed"
$ws.Range("C16").Value = "Open"

# Match the updated selection state from the diff
$ws.Range("C17").Select()
